$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The two single-cell "unchanged generated file" rows (01_main_user_schema_sequence.sql
# and 02_db_config_tools.sql) are templates for files that are generated but never
# change, so their template rows are removed entirely. Deleting row 34 twice removes
# both rows (the second row shifts up into position 34 after the first delete).
$ws.Rows.Item(34).Delete()
$ws.Rows.Item(34).Delete()

# The header row for the detail table grows a touch taller in the edited workbook.
$ws.Rows.Item(33).RowHeight = 16.5

# Reflect the new selection/active cell left after removing those rows.
[void]$ws.Range("A34:XFD35").Select()
